$d = $word.ActiveDocument

# 1) "have stake real money" -> "have staked real money"
$d.Content.Find.Execute("have stake real money", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "have staked real money", 2)

# 2) "something similar." -> "something similar that would appeal to the main player base."
$d.Content.Find.Execute("such as a tiny bonus in game or something similar.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "such as a tiny bonus in game or something similar that would appeal to the main player base.", 2)
